$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The run added averaged-intensity results for three new spiral sampling
# schemes (Gaussian-Quadrature keeps its data but now sorts ahead of the
# "NoRotation/Rotation" block; three new "Spiral-..." rows are inserted right
# after it). Insert 3 blank rows after row 10 so rows 11-16 (Rotation/HexGrid)
# shift down to rows 14-19.
$ws.Rows("11:13").Insert()

# Inherit the bold/bordered "index" style (col A, and col B uses default) from
# row 10 onto the freshly inserted rows 11-13.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$cols = @("C","D","E","F","G","H","I","J","K","L","M")

# row => @(index, label, C..M values)
$rows = @(
    @(10, 8, "Gaussian-Quadrature", @(1.003347938170875, 0.9733256173515474, 1.001125728766453, 1.003347938170875, 0.9836962612942944, 1.012687467361678, 0.9998331504704113, 0.9733256173515474, 0.9872256730590003, 0.9952868056149375, 0.9956693605692098)),
    @(11, 9, "Spiral-90deg-10rot-5space", @(1.001932495920071, 0.9649877851129955, 1.005231144843666, 1.001932495920071, 0.9792674103079102, 1.02038144023133, 1.003790213844324, 0.9649877851129955, 0.9851094649783307, 0.9935209804492009, 0.9959317483767162)),
    @(12, 10, "Spiral-90deg-15rot-5space", @(1.002242062709764, 0.9637150356918213, 1.005548370810161, 1.002242062709764, 0.978575604970431, 1.021194504960613, 1.004107201290522, 0.9637150356918213, 0.984631703250991, 0.9934368829803772, 0.9958971300722186)),
    @(13, 11, "Spiral-90deg-10rot-3space", @(1.001980459520982, 0.9646912219784327, 1.005332058075448, 1.001980459520982, 0.9791013684398893, 1.020561968197788, 1.00386947620977, 0.9646912219784327, 0.9850116400269404, 0.9934960497739611, 0.9959227587370517)),
    @(14, 12, "NoRotation-tilt60deg", @(1.045612, 0.8486359999999988, 1.026012, 1.045612, 0.9185960000000004, 1.089927999999998, 1.029976000000001, 0.8486359999999988, 0.9373239999999996, 0.9914679999999998, 0.9931266666666664)),
    @(15, 13, "Rotation-NoTilt", @(1.05, 0.88, 1.01, 1.05, 0.94, 1.06, 1.02, 0.88, 0.9450000000000001, 0.9975, 0.9933333333333335)),
    @(16, 14, "Rotation-60detTilt", @(1.026294226534397, 0.9282446038016006, 1.004314014310392, 1.026294226534397, 0.9632093075456044, 1.032200858931205, 1.009603000729594, 0.9282446038016006, 0.9662793090559965, 0.9962867677951965, 0.9939776686421323)),
    @(17, 15, "HexGrid-90degTilt5degRes", @(0.9946101379994868, 0.9953088850461074, 0.9952732980393554, 0.9946101379994868, 0.995222329733299, 0.9953607396813189, 0.9953620408711823, 0.9953088850461074, 0.9952910915427313, 0.9949506147711091, 0.9951895718951249)),
    @(18, 16, "HexGrid-90degTilt22p5degRes", @(0.9961260033708847, 0.9970693168684011, 0.9945301844045998, 0.9961260033708847, 0.9963751871146853, 0.9923561936083879, 0.9945017866503516, 0.9970693168684011, 0.9957997506365004, 0.9959628770036926, 0.9951597786695517)),
    @(19, 17, "HexGrid-60degTilt5degRes", @(0.9922145835281826, 1.006329311554197, 0.9926450763660348, 0.9922145835281826, 1.001982252353012, 0.9876881727365987, 0.9920151109149762, 1.006329311554197, 0.9994871939601159, 0.9958508887441493, 0.9954790845755003))
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $idx = $entry[1]
    $label = $entry[2]
    $vals = $entry[3]
    $ws.Range("A" + $r).Value = $idx
    $ws.Range("B" + $r).Value = $label
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}